$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 2371.7144
$ws.Range("I32").Value = 1667.3334
$ws.Range("J32").Value = 2900
$ws.Range("K32").Value = 1667.3334
$ws.Range("L32").Value = 2900
$ws.Range("M32").Value = -1341.3334
$ws.Range("N32").Value = -3552
$ws.Range("H33").Value = 692.0968
$ws.Range("I33").Value = 871.05554
$ws.Range("J33").Value = 444.30768
$ws.Range("K33").Value = 871.05554
$ws.Range("L33").Value = 444.30768
$ws.Range("M33").Value = -642.05554
$ws.Range("N33").Value = -902.30768
$ws.Range("H51").Value = 30311862
$ws.Range("I51").Value = 45457790
$ws.Range("J51").Value = 20000
$ws.Range("K51").Value = 45457790
$ws.Range("L51").Value = 20000
$ws.Range("M51").Value = -45457306
$ws.Range("N51").Value = -20968
$ws.Range("H58").Value = 2575.3333
$ws.Range("I58").Value = 250
$ws.Range("J58").Value = 2712.1177
$ws.Range("K58").Value = 750
$ws.Range("L58").Value = 8136.353099999999
$ws.Range("M58").Value = -600
$ws.Range("N58").Value = -8436.3531
$ws.Range("H62").Value = 18583.955
$ws.Range("I62").Value = 5207.0713
$ws.Range("J62").Value = 41993.5
$ws.Range("K62").Value = 5207.0713
$ws.Range("L62").Value = 41993.5
$ws.Range("M62").Value = -4583.0713
$ws.Range("N62").Value = -43241.5
$ws.Range("H65").Value = 18583.955
$ws.Range("I65").Value = 5207.0713
$ws.Range("J65").Value = 41993.5
$ws.Range("K65").Value = 26035.3565
$ws.Range("L65").Value = 209967.5
$ws.Range("M65").Value = -22915.3565
$ws.Range("N65").Value = -216207.5
$ws.Range("H116").Value = 7857.143
$ws.Range("I116").Value = 4600
$ws.Range("J116").Value = 16000
$ws.Range("K116").Value = 4600
$ws.Range("L116").Value = 16000
$ws.Range("M116").Value = -1158
$ws.Range("N116").Value = -22884
$ws.Range("H137").Value = 1837.2609
$ws.Range("I137").Value = 2126.5
$ws.Range("J137").Value = 1572.125
$ws.Range("K137").Value = 6379.5
$ws.Range("L137").Value = 4716.375
$ws.Range("M137").Value = -3829.5
$ws.Range("N137").Value = -9816.375
$ws.Range("H138").Value = 1857.3522
$ws.Range("I138").Value = 902.42224
$ws.Range("J138").Value = 3510.1155
$ws.Range("K138").Value = 2707.26672
$ws.Range("L138").Value = 10530.3465
$ws.Range("M138").Value = 2432.73328
$ws.Range("N138").Value = -20810.3465

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 499.5
$ws.Range("I3").Value = 499.5
$ws.Range("K3").Value = 499.5
$ws.Range("M3").Value = -384.5
$ws.Range("H124").Value = 20182.9
$ws.Range("J124").Value = 20182.9
$ws.Range("L124").Value = 20182.9
$ws.Range("N124").Value = -30002.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 2701
$ws.Range("I8").Value = 334.66666
$ws.Range("J8").Value = 9800
$ws.Range("K8").Value = 334.66666
$ws.Range("L8").Value = 9800
$ws.Range("M8").Value = -194.66666
$ws.Range("N8").Value = -10080

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9806086
$ws.Range("I31").Value = 1555.2593
$ws.Range("J31").Value = 47623560
$ws.Range("K31").Value = 1555.2593
$ws.Range("L31").Value = 47623560
$ws.Range("M31").Value = -1260.2593
$ws.Range("N31").Value = -47624150
$ws.Range("H34").Value = 9806086
$ws.Range("I34").Value = 1555.2593
$ws.Range("J34").Value = 47623560
$ws.Range("K34").Value = 1555.2593
$ws.Range("L34").Value = 47623560
$ws.Range("M34").Value = -1353.2593
$ws.Range("N34").Value = -47623964
$ws.Range("H50").Value = 10000
$ws.Range("I50").Value = 10000
$ws.Range("J50").Value = 0
$ws.Range("K50").Value = 10000
$ws.Range("L50").Value = 0
$ws.Range("M50").Value = -9375
$ws.Range("N50").Value = $null
$ws.Range("H62").Value = 19542
$ws.Range("I62").Value = 18824.334
$ws.Range("J62").Value = 20157.143
$ws.Range("K62").Value = 18824.334
$ws.Range("L62").Value = 20157.143
$ws.Range("M62").Value = -18200.334
$ws.Range("N62").Value = -21405.143
$ws.Range("H65").Value = 19542
$ws.Range("I65").Value = 18824.334
$ws.Range("J65").Value = 20157.143
$ws.Range("K65").Value = 94121.67
$ws.Range("L65").Value = 100785.715
$ws.Range("M65").Value = -91001.67
$ws.Range("N65").Value = -107025.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2294.6667
$ws.Range("I5").Value = 1004
$ws.Range("J5").Value = 2940
$ws.Range("K5").Value = 3012
$ws.Range("L5").Value = 8820
$ws.Range("M5").Value = -2900
$ws.Range("N5").Value = -9044
$ws.Range("H50").Value = 1285.6154
$ws.Range("I50").Value = 70.375
$ws.Range("K50").Value = 211.125
$ws.Range("M50").Value = 269.875
$ws.Range("H53").Value = 1285.6154
$ws.Range("I53").Value = 70.375
$ws.Range("K53").Value = 211.125
$ws.Range("M53").Value = 269.875
$ws.Range("H61").Value = 262.16666
$ws.Range("I61").Value = 262.16666
$ws.Range("K61").Value = 786.4999799999999
$ws.Range("M61").Value = -571.4999799999999
$ws.Range("H63").Value = 4816.4287
$ws.Range("I63").Value = 7101
$ws.Range("J63").Value = 4435.6665
$ws.Range("K63").Value = 21303
$ws.Range("L63").Value = 13306.9995
$ws.Range("M63").Value = -20554
$ws.Range("N63").Value = -14804.9995
$ws.Range("H64").Value = 3010
$ws.Range("I64").Value = 0
$ws.Range("J64").Value = 3010
$ws.Range("K64").Value = 0
$ws.Range("L64").Value = 9030
$ws.Range("M64").Value = $null
$ws.Range("N64").Value = -9570
$ws.Range("H66").Value = 4816.4287
$ws.Range("I66").Value = 7101
$ws.Range("J66").Value = 4435.6665
$ws.Range("K66").Value = 63909
$ws.Range("L66").Value = 39920.9985
$ws.Range("M66").Value = -60165
$ws.Range("N66").Value = -47408.9985
$ws.Range("H67").Value = 3010
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3010
$ws.Range("K67").Value = 0
$ws.Range("L67").Value = 9030
$ws.Range("M67").Value = $null
$ws.Range("N67").Value = -10902
$ws.Range("H69").Value = 2544.111
$ws.Range("I69").Value = 1345.75
$ws.Range("J69").Value = 3502.8
$ws.Range("K69").Value = 4037.25
$ws.Range("L69").Value = 10508.4
$ws.Range("M69").Value = -3226.25
$ws.Range("N69").Value = -12130.4
$ws.Range("H70").Value = 2186.5715
$ws.Range("I70").Value = 1173.1428
$ws.Range("J70").Value = 3200
$ws.Range("K70").Value = 3519.4284
$ws.Range("L70").Value = 9600
$ws.Range("M70").Value = -3204.4284
$ws.Range("N70").Value = -10230
$ws.Range("H72").Value = 2544.111
$ws.Range("I72").Value = 1345.75
$ws.Range("J72").Value = 3502.8
$ws.Range("K72").Value = 12111.75
$ws.Range("L72").Value = 31525.2
$ws.Range("M72").Value = -8055.75
$ws.Range("N72").Value = -39637.2
$ws.Range("H73").Value = 2186.5715
$ws.Range("I73").Value = 1173.1428
$ws.Range("J73").Value = 3200
$ws.Range("K73").Value = 3519.4284
$ws.Range("L73").Value = 9600
$ws.Range("M73").Value = -2427.4284
$ws.Range("N73").Value = -11784
$ws.Range("H75").Value = 3545.3635
$ws.Range("I75").Value = 866.3333
$ws.Range("J75").Value = 4550
$ws.Range("K75").Value = 2598.9999
$ws.Range("L75").Value = 13650
$ws.Range("M75").Value = -1600.9999
$ws.Range("N75").Value = -15646
$ws.Range("H76").Value = 3238.3333
$ws.Range("J76").Value = 3433.3333
$ws.Range("L76").Value = 10299.9999
$ws.Range("N76").Value = -11065.9999
$ws.Range("H78").Value = 3545.3635
$ws.Range("I78").Value = 866.3333
$ws.Range("J78").Value = 4550
$ws.Range("K78").Value = 7796.9997
$ws.Range("L78").Value = 40950
$ws.Range("M78").Value = -2804.9997
$ws.Range("N78").Value = -50934
$ws.Range("H79").Value = 3238.3333
$ws.Range("J79").Value = 3433.3333
$ws.Range("L79").Value = 10299.9999
$ws.Range("N79").Value = -12951.9999
$ws.Range("H87").Value = 4302.3335
$ws.Range("I87").Value = 925.6
$ws.Range("J87").Value = 6714.2856
$ws.Range("K87").Value = 2776.8
$ws.Range("L87").Value = 20142.8568
$ws.Range("M87").Value = -1528.8
$ws.Range("N87").Value = -22638.8568
$ws.Range("H88").Value = 4091.6667
$ws.Range("J88").Value = 4091.6667
$ws.Range("L88").Value = 12275.0001
$ws.Range("N88").Value = -13131.0001
$ws.Range("H90").Value = 4302.3335
$ws.Range("I90").Value = 925.6
$ws.Range("J90").Value = 6714.2856
$ws.Range("K90").Value = 8330.4
$ws.Range("L90").Value = 60428.5704
$ws.Range("M90").Value = -2090.4
$ws.Range("N90").Value = -72908.5704
$ws.Range("H91").Value = 4091.6667
$ws.Range("J91").Value = 4091.6667
$ws.Range("L91").Value = 12275.0001
$ws.Range("N91").Value = -15239.0001
$ws.Range("H122").Value = 2981.9285
$ws.Range("I122").Value = 340
$ws.Range("J122").Value = 3556.261
$ws.Range("K122").Value = 3060
$ws.Range("L122").Value = 32006.349
$ws.Range("M122").Value = -610
$ws.Range("N122").Value = -36906.349
$ws.Range("H135").Value = 2294.6667
$ws.Range("I135").Value = 1004
$ws.Range("J135").Value = 2940
$ws.Range("K135").Value = 9036
$ws.Range("L135").Value = 26460
$ws.Range("M135").Value = -6501
$ws.Range("N135").Value = -31530

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("N127").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 28596.666
$ws.Range("J125").Value = 28596.666
$ws.Range("L125").Value = 28596.666
$ws.Range("N125").Value = -38436.666
